$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking Price values so they remain stored as text
# (matching the source data which is t="inlineStr"/shared-string text, not numbers).
foreach ($cell in @("D5", "D6", "D8", "D11", "D13", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D28", "D30", "D32", "D33", "D35", "D38", "D39", "D43", "D44", "D45", "D50")) {
    $ws.Range($cell).NumberFormat = "@"
}

# Apply the updated cell values from the crypto data refresh
$ws.Range("D2").Value = "69.402.44"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "3.434.25"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "583.83"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").Value = "179.48"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "0.594"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  +7.10%  "
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("D11").Value = "48.59"
$ws.Range("E12").Value = "  +2.01%  "
$ws.Range("D13").Value = "686.51"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").Value = "3.986.38"
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("E15").Value = "  +1.73%  "
$ws.Range("D16").Value = "69.487.37"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "3.435.27"
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").Value = "17.89"
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("D20").Value = "11.35"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").Value = "0.913"
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("D22").Value = "5.38"
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("D23").Value = "17.06"
$ws.Range("D24").Value = "101.36"
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").Value = "3.92"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("D28").Value = "33.67"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("E29").Value = "  +2.01%  "
$ws.Range("D30").Value = "6.91"
$ws.Range("E30").Value = "  -1.93%  "
$ws.Range("E31").Value = "  +6.67%  "
$ws.Range("D32").Value = "562.48"
$ws.Range("E32").Value = "  +1.43%  "
$ws.Range("D33").Value = "11.06"
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("E34").Value = "  -0.71%  "
$ws.Range("D35").Value = "58.24"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").Value = "3.627.06"
$ws.Range("E37").Value = "  -2.55%  "
$ws.Range("D38").Value = "0.140"
$ws.Range("E38").Value = "  -2.13%  "
$ws.Range("D39").Value = "35.24"
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("D40").Value = "0.0₃0746"
$ws.Range("E40").Value = "  +9.00%  "
$ws.Range("E41").Value = "  +2.58%  "
$ws.Range("E42").Value = "  +1.95%  "
$ws.Range("D43").Value = "0.0425"
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "3.34"
$ws.Range("E44").Value = "  +2.99%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "0.336"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("E48").Value = "  +4.45%  "
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").Value = "131.38"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("E51").Value = "  +2.06%  "
